# Updates cryptos list values (Price / Volume(1h)) and the ONDO/SuiNetwork row swap
# to match the refreshed data, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.418.58"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").Value = "3.693.34"
$ws.Range("E3").Value = "  -3.46%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'686.47"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("D6").Value = "'161.16"
$ws.Range("E6").Value = "  -6.03%  "
$ws.Range("D7").Value = "3.691.94"
$ws.Range("E7").Value = "  -3.51%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -5.98%  "
$ws.Range("E10").Value = "  -8.49%  "
$ws.Range("D11").Value = "'7.24"
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("E12").Value = "  -10.29%  "
$ws.Range("D13").Value = "'0.0000235"
$ws.Range("E13").Value = "  -6.90%  "
$ws.Range("D14").Value = "4.317.58"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").Value = "'32.66"
$ws.Range("E15").Value = "  -11.34%  "
$ws.Range("D16").Value = "3.688.66"
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("D17").Value = "69.429.36"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("E19").Value = "  -9.76%  "
$ws.Range("D20").Value = "'6.47"
$ws.Range("E20").Value = "  -11.00%  "
$ws.Range("D21").Value = "'474.91"
$ws.Range("E21").Value = "  -7.69%  "
$ws.Range("D22").Value = "'9.91"
$ws.Range("E22").Value = "  -6.15%  "
$ws.Range("D23").Value = "'0.649"
$ws.Range("E23").Value = "  -9.70%  "
$ws.Range("D24").Value = "'79.67"
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("D25").Value = "3.836.51"
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'0.0000126"
$ws.Range("E27").Value = "  -11.23%  "
$ws.Range("D28").Value = "'11.06"
$ws.Range("E28").Value = "  -13.60%  "
$ws.Range("D29").Value = "'9.25"
$ws.Range("E29").Value = "  -11.01%  "
$ws.Range("E30").Value = "  -9.85%  "
$ws.Range("D31").Value = "'1.77"
$ws.Range("E31").Value = "  -11.88%  "
$ws.Range("D32").Value = "'6.70"
$ws.Range("E32").Value = "  -9.55%  "
$ws.Range("E33").Value = "  -10.81%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'26.79"
$ws.Range("E35").Value = "  -8.87%  "
$ws.Range("E36").Value = "  -7.12%  "
$ws.Range("D37").Value = "'8.25"
$ws.Range("E37").Value = "  -12.05%  "
$ws.Range("D38").Value = "'6.13"
$ws.Range("E38").Value = "  -8.22%  "
$ws.Range("D39").Value = "'2.29"
$ws.Range("E39").Value = "  -5.03%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'0.0911"
$ws.Range("E41").Value = "  -10.32%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'0.944"
$ws.Range("E43").Value = "  -6.94%  "
$ws.Range("D44").Value = "'166.34"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'47.93"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("D46").Value = "'2.75"
$ws.Range("E46").Value = "  -15.56%  "
$ws.Range("B47").Value = "SuiNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D47").Value = "'1.12"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.32"
$ws.Range("E48").Value = "  -4.76%  "
$ws.Range("D49").Value = "'28.57"
$ws.Range("E49").Value = "  -7.56%  "
$ws.Range("D50").Value = "'0.000277"
$ws.Range("E50").Value = "  -8.95%  "
$ws.Range("D51").Value = "'7.86"
$ws.Range("E51").Value = "  -9.50%  "
